$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2488
$ws.Cells.Item(137, 9).Value = 1735.8823
$ws.Cells.Item(137, 10).Value = 4314.5713
$ws.Cells.Item(137, 11).Value = 5207.6469
$ws.Cells.Item(137, 12).Value = 12943.7139
$ws.Cells.Item(137, 13).Value = -2657.6469
$ws.Cells.Item(137, 14).Value = -18043.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 2950.6667
$ws.Cells.Item(6, 9).Value = 601
$ws.Cells.Item(6, 10).Value = 7650
$ws.Cells.Item(6, 11).Value = 601
$ws.Cells.Item(6, 12).Value = 7650
$ws.Cells.Item(6, 13).Value = -428
$ws.Cells.Item(6, 14).Value = -7996
$ws.Cells.Item(32, 8).Value = 6413.8735
$ws.Cells.Item(32, 9).Value = 5707.232
$ws.Cells.Item(32, 10).Value = 18002.8
$ws.Cells.Item(32, 11).Value = 5707.232
$ws.Cells.Item(32, 12).Value = 18002.8
$ws.Cells.Item(32, 13).Value = -5420.232
$ws.Cells.Item(32, 14).Value = -18576.8
$ws.Cells.Item(61, 8).Value = 5113.82
$ws.Cells.Item(61, 9).Value = 9223.4
$ws.Cells.Item(61, 10).Value = 3352.5715
$ws.Cells.Item(61, 11).Value = 9223.4
$ws.Cells.Item(61, 12).Value = 3352.5715
$ws.Cells.Item(61, 13).Value = -9011.4
$ws.Cells.Item(61, 14).Value = -3776.5715
$ws.Cells.Item(74, 8).Value = 1555.9762
$ws.Cells.Item(74, 9).Value = 1231.7587
$ws.Cells.Item(74, 11).Value = 1231.7587
$ws.Cells.Item(74, 13).Value = -357.7587000000001
$ws.Cells.Item(77, 8).Value = 1555.9762
$ws.Cells.Item(77, 9).Value = 1231.7587
$ws.Cells.Item(77, 11).Value = 6158.793500000001
$ws.Cells.Item(77, 13).Value = -1790.793500000001
$ws.Cells.Item(132, 8).Value = 2657.9312
$ws.Cells.Item(132, 9).Value = 1684.963
$ws.Cells.Item(132, 10).Value = 3505.3547
$ws.Cells.Item(132, 11).Value = 5054.889
$ws.Cells.Item(132, 12).Value = 10516.0641
$ws.Cells.Item(132, 13).Value = -2524.889
$ws.Cells.Item(132, 14).Value = -15576.0641
$ws.Cells.Item(136, 8).Value = 5113.82
$ws.Cells.Item(136, 9).Value = 9223.4
$ws.Cells.Item(136, 10).Value = 3352.5715
$ws.Cells.Item(136, 11).Value = 27670.2
$ws.Cells.Item(136, 12).Value = 10057.7145
$ws.Cells.Item(136, 13).Value = -25120.2
$ws.Cells.Item(136, 14).Value = -15157.7145
$ws.Cells.Item(139, 8).Value = 38569.215
$ws.Cells.Item(139, 10).Value = 38569.215
$ws.Cells.Item(139, 12).Value = 38569.215
$ws.Cells.Item(139, 14).Value = -48849.215

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 9914.25
$ws.Cells.Item(105, 9).Value = 14776.5625
$ws.Cells.Item(105, 11).Value = 14776.5625
$ws.Cells.Item(105, 13).Value = -13029.5625
$ws.Cells.Item(134, 8).Value = 4172.7393
$ws.Cells.Item(134, 9).Value = 4854.4
$ws.Cells.Item(134, 10).Value = 2894.625
$ws.Cells.Item(134, 11).Value = 14563.2
$ws.Cells.Item(134, 12).Value = 8683.875
$ws.Cells.Item(134, 13).Value = -12028.2
$ws.Cells.Item(134, 14).Value = -13753.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2886.1526
$ws.Cells.Item(31, 9).Value = 2136.9033
$ws.Cells.Item(31, 10).Value = 3715.6785
$ws.Cells.Item(31, 11).Value = 2136.9033
$ws.Cells.Item(31, 12).Value = 3715.6785
$ws.Cells.Item(31, 13).Value = -1841.9033
$ws.Cells.Item(31, 14).Value = -4305.6785
$ws.Cells.Item(34, 8).Value = 2886.1526
$ws.Cells.Item(34, 9).Value = 2136.9033
$ws.Cells.Item(34, 10).Value = 3715.6785
$ws.Cells.Item(34, 11).Value = 2136.9033
$ws.Cells.Item(34, 12).Value = 3715.6785
$ws.Cells.Item(34, 13).Value = -1934.9033
$ws.Cells.Item(34, 14).Value = -4119.6785
$ws.Cells.Item(58, 8).Value = 2077.257
$ws.Cells.Item(58, 9).Value = 1869.2354
$ws.Cells.Item(58, 10).Value = 2273.7222
$ws.Cells.Item(58, 11).Value = 1869.2354
$ws.Cells.Item(58, 12).Value = 2273.7222
$ws.Cells.Item(58, 13).Value = -1666.2354
$ws.Cells.Item(58, 14).Value = -2679.7222
$ws.Cells.Item(105, 8).Value = 2464.5
$ws.Cells.Item(105, 9).Value = 2329.8572
$ws.Cells.Item(105, 11).Value = 2329.8572
$ws.Cells.Item(105, 13).Value = -582.8571999999999
$ws.Cells.Item(132, 8).Value = 1969.2354
$ws.Cells.Item(132, 9).Value = 1531.6316
$ws.Cells.Item(132, 10).Value = 2523.5334
$ws.Cells.Item(132, 11).Value = 4594.8948
$ws.Cells.Item(132, 12).Value = 7570.600199999999
$ws.Cells.Item(132, 13).Value = -2064.8948
$ws.Cells.Item(132, 14).Value = -12630.6002
$ws.Cells.Item(134, 8).Value = 3244.7715
$ws.Cells.Item(134, 9).Value = 3419.7307
$ws.Cells.Item(134, 10).Value = 2739.3333
$ws.Cells.Item(134, 11).Value = 10259.1921
$ws.Cells.Item(134, 12).Value = 8217.999899999999
$ws.Cells.Item(134, 13).Value = -7724.1921
$ws.Cells.Item(134, 14).Value = -13287.9999
$ws.Cells.Item(136, 8).Value = 2077.257
$ws.Cells.Item(136, 9).Value = 1869.2354
$ws.Cells.Item(136, 10).Value = 2273.7222
$ws.Cells.Item(136, 11).Value = 5607.706200000001
$ws.Cells.Item(136, 12).Value = 6821.1666
$ws.Cells.Item(136, 13).Value = -3057.706200000001
$ws.Cells.Item(136, 14).Value = -11921.1666
$ws.Cells.Item(138, 8).Value = 52840
$ws.Cells.Item(138, 10).Value = 52840
$ws.Cells.Item(138, 12).Value = 52840
$ws.Cells.Item(138, 14).Value = -63120

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 2116
$ws.Cells.Item(69, 10).Value = 2275.4285
$ws.Cells.Item(69, 12).Value = 6826.2855
$ws.Cells.Item(69, 14).Value = -8448.2855
$ws.Cells.Item(72, 8).Value = 2116
$ws.Cells.Item(72, 10).Value = 2275.4285
$ws.Cells.Item(72, 12).Value = 20478.8565
$ws.Cells.Item(72, 14).Value = -28590.8565
$ws.Cells.Item(113, 8).Value = 1053206.5
$ws.Cells.Item(113, 9).Value = 1163383.5
$ws.Cells.Item(113, 10).Value = 714805.6
$ws.Cells.Item(113, 11).Value = 3490150.5
$ws.Cells.Item(113, 12).Value = 2144416.8
$ws.Cells.Item(113, 13).Value = -3487980.5
$ws.Cells.Item(113, 14).Value = -2148756.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(98, 8).Value = 40728.6
$ws.Cells.Item(98, 10).Value = 40728.6
$ws.Cells.Item(98, 12).Value = 40728.6
$ws.Cells.Item(98, 14).Value = -46718.6
$ws.Cells.Item(102, 8).Value = 893546.3
$ws.Cells.Item(102, 9).Value = 1304760
$ws.Cells.Item(102, 10).Value = 2583.3333
$ws.Cells.Item(102, 11).Value = 1304760
$ws.Cells.Item(102, 12).Value = 2583.3333
$ws.Cells.Item(102, 13).Value = -1303138
$ws.Cells.Item(102, 14).Value = -5827.3333
$ws.Cells.Item(105, 8).Value = 30970
$ws.Cells.Item(105, 10).Value = 30970
$ws.Cells.Item(105, 12).Value = 30970
$ws.Cells.Item(105, 14).Value = -37958
$ws.Cells.Item(122, 8).Value = 13230194
$ws.Cells.Item(122, 9).Value = 1737898.5
$ws.Cells.Item(122, 10).Value = 50005540
$ws.Cells.Item(122, 11).Value = 5213695.5
$ws.Cells.Item(122, 12).Value = 150016620
$ws.Cells.Item(122, 13).Value = -5211245.5
$ws.Cells.Item(122, 14).Value = -150021520
$ws.Cells.Item(123, 8).Value = 16262.387
$ws.Cells.Item(123, 10).Value = 16431.28
$ws.Cells.Item(123, 12).Value = 16431.28
$ws.Cells.Item(123, 14).Value = -21331.28
$ws.Cells.Item(126, 8).Value = 8688.25
$ws.Cells.Item(126, 9).Value = 11082.909
$ws.Cells.Item(126, 10).Value = 3420
$ws.Cells.Item(126, 11).Value = 33248.727
$ws.Cells.Item(126, 12).Value = 10260
$ws.Cells.Item(126, 13).Value = -30778.727
$ws.Cells.Item(126, 14).Value = -15200
$ws.Cells.Item(132, 8).Value = 30066
$ws.Cells.Item(132, 9).Value = 62353.06
$ws.Cells.Item(132, 10).Value = 2622
$ws.Cells.Item(132, 11).Value = 187059.18
$ws.Cells.Item(132, 12).Value = 7866
$ws.Cells.Item(132, 13).Value = -184529.18
$ws.Cells.Item(132, 14).Value = -12926

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 28573362
$ws.Cells.Item(40, 9).Value = 43480340
$ws.Cells.Item(40, 10).Value = 1650.4166
$ws.Cells.Item(40, 11).Value = 43480340
$ws.Cells.Item(40, 12).Value = 1650.4166
$ws.Cells.Item(40, 13).Value = -43480204
$ws.Cells.Item(40, 14).Value = -1922.4166
$ws.Cells.Item(46, 8).Value = 1437.75
$ws.Cells.Item(46, 9).Value = 1300
$ws.Cells.Item(46, 10).Value = 1483.6666
$ws.Cells.Item(46, 11).Value = 1300
$ws.Cells.Item(46, 12).Value = 1483.6666
$ws.Cells.Item(46, 13).Value = -1112
$ws.Cells.Item(46, 14).Value = -1859.6666
$ws.Cells.Item(132, 8).Value = 10758385
$ws.Cells.Item(132, 9).Value = 18526816
$ws.Cells.Item(132, 10).Value = 2096.8462
$ws.Cells.Item(132, 11).Value = 55580448
$ws.Cells.Item(132, 12).Value = 6290.5386
$ws.Cells.Item(132, 13).Value = -55577918
$ws.Cells.Item(132, 14).Value = -11350.5386
$ws.Cells.Item(139, 8).Value = 45678.75
$ws.Cells.Item(139, 10).Value = 45678.75
$ws.Cells.Item(139, 12).Value = 45678.75
$ws.Cells.Item(139, 14).Value = -55958.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 26125
$ws.Cells.Item(41, 10).Value = 26125
$ws.Cells.Item(41, 12).Value = 26125
$ws.Cells.Item(41, 14).Value = -26905
$ws.Cells.Item(122, 8).Value = 1914.8462
$ws.Cells.Item(122, 9).Value = 1288.3
$ws.Cells.Item(122, 10).Value = 4003.3333
$ws.Cells.Item(122, 11).Value = 3864.9
$ws.Cells.Item(122, 12).Value = 12009.9999
$ws.Cells.Item(122, 13).Value = -1414.9
$ws.Cells.Item(122, 14).Value = -16909.9999
$ws.Cells.Item(133, 8).Value = 40642.168
$ws.Cells.Item(133, 10).Value = 40642.168
$ws.Cells.Item(133, 12).Value = 40642.168
$ws.Cells.Item(133, 14).Value = -50762.168
$ws.Cells.Item(135, 8).Value = 40715
$ws.Cells.Item(135, 10).Value = 40715
$ws.Cells.Item(135, 12).Value = 40715
$ws.Cells.Item(135, 14).Value = -50855
$ws.Cells.Item(141, 8).Value = 57466.875
$ws.Cells.Item(141, 10).Value = 57466.875
$ws.Cells.Item(141, 12).Value = 57466.875
$ws.Cells.Item(141, 14).Value = -67826.875
